# Logged Week 16 and performed season sim from Week 17
#
# - Add the new player "M.Hall" as a fresh row on the WR sheet (row 8),
#   with this week's stats still at 0 (just logged, not yet simmed).
# - Move the active/selected sheet from RB to WR, with the selection
#   resting at G9 (just past the newly added row).

$wb = $excel.ActiveWorkbook

$wsWR = $wb.Worksheets.Item("WR")

# Log the new player on the WR sheet.
$wsWR.Cells.Item(8, 1).Value = "M.Hall"
for ($col = 2; $col -le 10; $col++) {
    $wsWR.Cells.Item(8, $col).Value = 0
}

# Season sim moved focus on to the WR tab, selecting the cell just below
# the newly logged row.
$wsWR.Select()
$wsWR.Range("G9").Select()
